$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be auto-parsed as numbers
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.908.00"
$ws.Range("E2").Value = "  +3.48%  "

$ws.Range("D3").Value = "3.265.27"
$ws.Range("E3").Value = "  +2.88%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "580.79"
$ws.Range("E5").Value = "  +1.68%  "

$ws.Range("D6").Value = "181.93"
$ws.Range("E6").Value = "  +6.03%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  +0.32%  "

$ws.Range("D9").Value = "3.266.08"
$ws.Range("E9").Value = "  +2.98%  "

$ws.Range("E10").Value = "  +7.10%  "

$ws.Range("E11").Value = "  +2.76%  "

$ws.Range("D12").Value = "0.417"
$ws.Range("E12").Value = "  +6.48%  "

$ws.Range("D13").Value = "3.833.85"
$ws.Range("E13").Value = "  +3.06%  "

$ws.Range("D14").Value = "0.138"
$ws.Range("E14").Value = "  +1.04%  "

$ws.Range("D15").Value = "28.46"
$ws.Range("E15").Value = "  +4.24%  "

$ws.Range("D16").Value = "67.875.44"
$ws.Range("E16").Value = "  +3.58%  "

$ws.Range("E17").Value = "  +3.50%  "

$ws.Range("D18").Value = "3.251.44"
$ws.Range("E18").Value = "  +2.39%  "

$ws.Range("D19").Value = "5.85"
$ws.Range("E19").Value = "  +2.54%  "

$ws.Range("D20").Value = "13.53"
$ws.Range("E20").Value = "  +5.21%  "

$ws.Range("D21").Value = "375.99"
$ws.Range("E21").Value = "  +4.81%  "

$ws.Range("D22").Value = "7.66"
$ws.Range("E22").Value = "  +5.24%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").Value = "71.17"

$ws.Range("E25").Value = "  +3.81%  "

$ws.Range("E26").Value = "  +4.43%  "

$ws.Range("E27").Value = "  -1.81%  "

$ws.Range("E28").Value = "  +2.28%  "

$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("E30").Value = "  +3.17%  "

$ws.Range("D31").Value = "5.69"
$ws.Range("E31").Value = "  +6.03%  "

$ws.Range("D32").Value = "22.80"
$ws.Range("E32").Value = "  +3.93%  "

$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("E34").Value = "  +6.06%  "

$ws.Range("D35").Value = "6.93"
$ws.Range("E35").Value = "  +4.91%  "

$ws.Range("E36").Value = "  +4.97%  "

$ws.Range("D37").Value = "161.70"
$ws.Range("E37").Value = "  +1.15%  "

$ws.Range("E38").Value = "  +2.32%  "

$ws.Range("D39").Value = "1.84"
$ws.Range("E39").Value = "  +3.19%  "

$ws.Range("D40").Value = "6.79"
$ws.Range("E40").Value = "  +11.15%  "

$ws.Range("D41").Value = "26.77"
$ws.Range("E41").Value = "  +1.37%  "

$ws.Range("E42").Value = "  +11.21%  "

$ws.Range("E43").Value = "  +4.59%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "25.78"
$ws.Range("E44").Value = "  +7.32%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.690.79"
$ws.Range("E45").Value = "  +1.65%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "350.77"
$ws.Range("E46").Value = "  +6.92%  "

$ws.Range("D47").Value = "40.81"
$ws.Range("E47").Value = "  +3.07%  "

$ws.Range("D48").Value = "0.0681"
$ws.Range("E48").Value = "  +3.62%  "

$ws.Range("E49").Value = "  +2.78%  "

$ws.Range("E50").Value = "  +5.97%  "

$ws.Range("D51").Value = "0.103"
$ws.Range("E51").Value = "  +0.52%  "
